$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2025-05-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-06 Tuesday", 2) | Out-Null

# Update each of the 100 math-expression cells in the table (row-major order,
# matches document order) directly via Cell.Range.Text so duplicate old values
# (e.g. two "11+88=99" cells with different replacements) resolve unambiguously.
$tbl = $d.Tables.Item(1)
$values = @(
    "64-21=43",
    "8+80=88",
    "21+33=54",
    "85+7=92",
    "87-1=86",
    "29+70=99",
    "26+44=70",
    "95-80=15",
    "31+12=43",
    "45+37=82",
    "21+12=33",
    "82-19=63",
    "91-9=82",
    "11-6=5",
    "20+0=20",
    "55-42=13",
    "46+49=95",
    "28+44=72",
    "56+20=76",
    "70-56=14",
    "53-21=32",
    "16+61=77",
    "69-11=58",
    "44+47=91",
    "41+33=74",
    "85-65=20",
    "41+2=43",
    "12+9=21",
    "41+5=46",
    "87-81=6",
    "83-12=71",
    "42-1=41",
    "35-18=17",
    "53-42=11",
    "58-11=47",
    "14+36=50",
    "56-24=32",
    "63-45=18",
    "66-38=28",
    "48-42=6",
    "95-34=61",
    "11+39=50",
    "31+14=45",
    "72+21=93",
    "11+10=21",
    "3+12=15",
    "32+18=50",
    "43-36=7",
    "91-55=36",
    "89-29=60",
    "61-28=33",
    "18-3=15",
    "93-54=39",
    "74-59=15",
    "13+33=46",
    "92-41=51",
    "97-64=33",
    "36+30=66",
    "39+49=88",
    "97-35=62",
    "37-6=31",
    "2+20=22",
    "98-91=7",
    "64-61=3",
    "58+33=91",
    "29+10=39",
    "79+8=87",
    "75+1=76",
    "14+2=16",
    "61+12=73",
    "90+2=92",
    "48+49=97",
    "72-42=30",
    "47-5=42",
    "3+19=22",
    "42+8=50",
    "9+78=87",
    "89+10=99",
    "4+54=58",
    "22-9=13",
    "37-7=30",
    "13+81=94",
    "10+63=73",
    "78-63=15",
    "63-35=28",
    "84-61=23",
    "14+35=49",
    "22-8=14",
    "76-37=39",
    "39+49=88",
    "76-26=50",
    "4+21=25",
    "0+43=43",
    "21-4=17",
    "59-6=53",
    "98-90=8",
    "83-32=51",
    "98-21=77",
    "94-81=13",
    "35+4=39"
)

$cols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [int][Math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $tbl.Cell($row, $col).Range.Text = $values[$i]
}

Write-Output "done"
